$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9, pushing the existing rows 9-11 down to 10-12
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new weekly entry (copy of row layout, new values)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44664
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112037
$ws.Range("G9").Value = "Cebollín"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 8500
$ws.Range("M9").Value = 8250
$ws.Range("N9").Value = "$/paquete 36 unidades"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 229
$ws.Range("Q9").Value = 36
$ws.Range("R9").Value = "Hortaliza"

# Match the style (date format) used in column D for the rest of the rows
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
